# Tripadvisor New Orleans shard workbook update
# - Reorder worksheet tabs so "review_info" precedes "hotel_info"
# - Insert a new "State" column into hotel_info (between Hotel_Name and City)
#   and populate it with "Louisiana" for the existing hotel row

$wb = $excel.ActiveWorkbook

$wsHotel = $wb.Worksheets.Item("hotel_info")
$wsReview = $wb.Worksheets.Item("review_info")

# Move review_info before hotel_info so the tab order becomes:
# review_info, hotel_info
$wsReview.Move($wsHotel)

# Worksheet object references can become stale (positional) after a Move,
# so re-acquire the hotel_info sheet by name before editing it further.
$wsHotel = $wb.Worksheets.Item("hotel_info")

# Insert a new blank column at column C (pushes City, Zip, ... one to the right)
$wsHotel.Columns.Item(3).Insert()

# Populate the new "State" column
$wsHotel.Cells.Item(1, 3).Value = "State"
$wsHotel.Cells.Item(2, 3).Value = "Louisiana"
